$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New Gini coefficient (column J) values for rows 5-10 (previously empty)
$ws1.Range("J5").Value = 0.427
$ws1.Range("J6").Value = 0.4323
$ws1.Range("J7").Value = 0.4337
$ws1.Range("J8").Value = 0.4318
$ws1.Range("J9").Value = 0.4377
$ws1.Range("J10").Value = 0.4201

# Updated Gini coefficient values for existing rows
$ws1.Range("J11").Value = 0.4381
$ws1.Range("J12").Value = 0.4274
$ws1.Range("J16").Value = 0.4284
$ws1.Range("J17").Value = 0.4294
$ws1.Range("J18").Value = 0.4313
$ws1.Range("J19").Value = 0.4321

# Sheet1 becomes the active/selected tab, with J13 as the active cell
$ws1.Activate()
$ws1.Range("J13").Select()
